$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

$data = @(
    @("L_Tooltip_1", "Подсказка #1", "Tooltip #1"),
    @("L_Tooltip_2", "Подсказка #2", "Tooltip #2"),
    @("L_Tooltip_3", "Подсказка #3", "Tooltip #3"),
    @("L_Tooltip_4", "Подсказка #4", "Tooltip #4"),
    @("L_Tooltip_5", "Подсказка #5", "Tooltip #5")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 5 + $i
    $ws.Range("A$r").Value = $data[$i][0]
    $ws.Range("B$r").Value = $data[$i][1]
    $ws.Range("C$r").Value = $data[$i][2]
}

$ws.Activate()
$ws.Range("D12").Select()

$char = $wb.Worksheets.Item("Character")
$char.Range("C10").Select()
